$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

function Get-WholeParagraphRange($findText) {
    $r = $d.Content.Duplicate
    $r.Start = $d.Content.Start
    $r.End = $d.Content.End
    $found = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return $null
    }
    $para = $r.Paragraphs(1)
    return $para.Range
}

# ---------------------------------------------------------------------------
# Change 1: paragraph "25" + ".0" -- drop the _GoBack bookmark that used to
# sit between the two runs.
# ---------------------------------------------------------------------------
$p1 = Get-WholeParagraphRange("25.0")
$rPr25 = "<w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/><w:b/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr>"
$frag1 = "<w:p $wns w:rsidR='00B20107' w:rsidRPr='00506885' w:rsidRDefault='0013348E' w:rsidP='00EA7375'>" +
         "<w:pPr><w:jc w:val='center'/>$rPr25</w:pPr>" +
         "<w:r>$rPr25<w:t>25</w:t></w:r>" +
         "<w:r w:rsidR='001C260E'>$rPr25<w:t>.0</w:t></w:r>" +
         "</w:p>"
[void]$p1.InsertXML($frag1)

# ---------------------------------------------------------------------------
# Change 2: "Alternative Flow of Events" heading paragraph -- the paragraph
# mark's run properties swap the themed/shaded blue for the flat accent1 blue.
# ---------------------------------------------------------------------------
$p2 = Get-WholeParagraphRange("Alternative Flow of Events")
$frag2 = "<w:p $wns w:rsidR='00605DF3' w:rsidRDefault='00B20107' w:rsidP='002165E1'>" +
         "<w:pPr><w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/><w:b/><w:color w:val='4F81BD' w:themeColor='accent1'/><w:sz w:val='24'/></w:rPr></w:pPr>" +
         "<w:r w:rsidRPr='00F115F3'><w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/><w:b/><w:color w:val='365F91' w:themeColor='accent1' w:themeShade='BF'/><w:sz w:val='24'/></w:rPr><w:t>Alternative Flow of Events</w:t></w:r>" +
         "</w:p>"
[void]$p2.InsertXML($frag2)

# ---------------------------------------------------------------------------
# Change 3: the old single "Line 1: User dies before beating the enemy"
# paragraph becomes three paragraphs: an empty bookmarked one (now hosting
# _GoBack), a bold "Sub Event" heading, and a plain "25.1: Throws the
# weapon." line.
# ---------------------------------------------------------------------------
$p3 = Get-WholeParagraphRange("Line 1: User dies before beating the enemy")
$rPrBold = "<w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/><w:b/><w:color w:val='365F91' w:themeColor='accent1' w:themeShade='BF'/><w:sz w:val='24'/></w:rPr>"
$frag3 = "<w:p $wns w:rsidR='007D6791' w:rsidRPr='007D6791' w:rsidRDefault='007D6791' w:rsidP='002165E1'>" +
         "<w:pPr>$rPrBold</w:pPr>" +
         "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
         "</w:p>" +
         "<w:p $wns><w:pPr>$rPrBold</w:pPr><w:r>$rPrBold<w:t>Sub Event</w:t></w:r></w:p>" +
         "<w:p $wns><w:pPr><w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/></w:rPr></w:pPr>" +
         "<w:r><w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/></w:rPr><w:t>25.1: Throws the weapon.</w:t></w:r></w:p>"
[void]$p3.InsertXML($frag3)

Write-Output "done"
